$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Shift the name/label/valueType (columns B:D) of rows 21..78 up into rows 20..77,
# effectively removing the 'ff10' / 'Total number of stillbirths' entry while leaving
# the index column (A) untouched, then drop the now-duplicated trailing row 78.
$ws.Range("B20").Value = "kind1"
$ws.Range("C20").Value = "How old were you when you gave birth to your first child?"
$ws.Range("D20").Value = "integer"
$ws.Range("B21").Value = "prevcanc"
$ws.Range("C21").Value = "prevalent cancer"
$ws.Range("D21").Value = "integer"
$ws.Range("B22").Value = "lipidlower"
$ws.Range("C22").Value = "Lipid-lowering drugs"
$ws.Range("D22").Value = "integer"
$ws.Range("B23").Value = "nsar_excl_ASS"
$ws.Range("C23").Value = "non-steroidal anti-inflammatory drug (excl.ASS)"
$ws.Range("D23").Value = "integer"
$ws.Range("B24").Value = "casemi_fup5"
$ws.Range("C24").Value = "case status of myocardial infarction at FUP5"
$ws.Range("D24").Value = "integer"
$ws.Range("B25").Value = "casestroke_fup5"
$ws.Range("C25").Value = "case status of stroke at FUP5"
$ws.Range("D25").Value = "integer"
$ws.Range("B26").Value = "dd_incmi"
$ws.Range("C26").Value = "date of diagnosis of myocardial infarction"
$ws.Range("D26").Value = "date"
$ws.Range("B27").Value = "dd_incstroke"
$ws.Range("C27").Value = "date of diagnosis of stroke"
$ws.Range("D27").Value = "date"
$ws.Range("B28").Value = "caseI63_fup5"
$ws.Range("C28").Value = "ischaemic insult (ICD-10 I63)"
$ws.Range("D28").Value = "integer"
$ws.Range("B29").Value = "caseI61_fup5"
$ws.Range("C29").Value = "intracerebral hemorrhage (ICD-10 I61)"
$ws.Range("D29").Value = "integer"
$ws.Range("B30").Value = "casehyp_fup5"
$ws.Range("C30").Value = "case status of essential hypertension at FUP5"
$ws.Range("D30").Value = "integer"
$ws.Range("B31").Value = "dd_inchyp"
$ws.Range("C31").Value = "date of diagnosis of essential hypertension"
$ws.Range("D31").Value = "date"
$ws.Range("B32").Value = "casehf_fup5"
$ws.Range("C32").Value = "case status of heart failure at FUP5"
$ws.Range("D32").Value = "integer"
$ws.Range("B33").Value = "dd_inchf"
$ws.Range("C33").Value = "date of diagnosis of heart failure"
$ws.Range("D33").Value = "date"
$ws.Range("B34").Value = "casediab_fup5"
$ws.Range("C34").Value = "case status of diabetes at FUP5"
$ws.Range("D34").Value = "integer"
$ws.Range("B35").Value = "dd_incdiab"
$ws.Range("C35").Value = "date of diagnosis of diabetes"
$ws.Range("D35").Value = "date"
$ws.Range("B36").Value = "inccanc_fup5"
$ws.Range("C36").Value = "incident first occuring cancer at FUP5"
$ws.Range("D36").Value = "integer"
$ws.Range("B37").Value = "dd_inccanc"
$ws.Range("C37").Value = "date of diagnosis of first occuring cancer"
$ws.Range("D37").Value = "date"
$ws.Range("B38").Value = "vitstat5"
$ws.Range("C38").Value = "vital status at FUP5"
$ws.Range("D38").Value = "integer"
$ws.Range("B39").Value = "age_death"
$ws.Range("C39").Value = "age of death"
$ws.Range("D39").Value = "decimal"
$ws.Range("B40").Value = "age_fup5"
$ws.Range("C40").Value = "age at FUP5"
$ws.Range("D40").Value = "decimal"
$ws.Range("B41").Value = "bmi0"
$ws.Range("C41").Value = "BMI at baseline"
$ws.Range("D41").Value = "decimal"
$ws.Range("B42").Value = "bmi_f4"
$ws.Range("C42").Value = "BMI at FUP4"
$ws.Range("D42").Value = "decimal"
$ws.Range("B43").Value = "waist0"
$ws.Range("C43").Value = "waist circumference at baseline [cm]"
$ws.Range("D43").Value = "decimal"
$ws.Range("B44").Value = "hip0"
$ws.Range("C44").Value = "hip circumference at baseline [cm]"
$ws.Range("D44").Value = "decimal"
$ws.Range("B45").Value = "waist_f4"
$ws.Range("C45").Value = "waist circumference at FUP4 [cm]"
$ws.Range("D45").Value = "decimal"
$ws.Range("B46").Value = "hip_f4"
$ws.Range("C46").Value = "hip circumference at FUP4 [cm]"
$ws.Range("D46").Value = "decimal"
$ws.Range("B47").Value = "age_anth_f4"
$ws.Range("C47").Value = "age of anthropometric measurement at FUP4"
$ws.Range("D47").Value = "decimal"
$ws.Range("B48").Value = "GJ"
$ws.Range("C48").Value = "Total energy intake at baseline [kJ/d]"
$ws.Range("D48").Value = "decimal"
$ws.Range("B49").Value = "corr_trigly"
$ws.Range("C49").Value = "corrected triglycerides [mg/dL]"
$ws.Range("D49").Value = "decimal"
$ws.Range("B50").Value = "corr_chol"
$ws.Range("C50").Value = "corrected cholesterol [mg/dL]"
$ws.Range("D50").Value = "decimal"
$ws.Range("B51").Value = "corr_hdl"
$ws.Range("C51").Value = "corrected HDL cholesterol [mg/dL]"
$ws.Range("D51").Value = "decimal"
$ws.Range("B52").Value = "ZK"
$ws.Range("C52").Value = "carbohydrate intake at baseline [g/d]"
$ws.Range("D52").Value = "decimal"
$ws.Range("B53").Value = "ZE"
$ws.Range("C53").Value = "protein intake at baseline [g/d]"
$ws.Range("D53").Value = "decimal"
$ws.Range("B54").Value = "ZF"
$ws.Range("C54").Value = "fat intake at baseline [g/d]"
$ws.Range("D54").Value = "decimal"
$ws.Range("B55").Value = "ZA"
$ws.Range("C55").Value = "alcohol intake at baseline [g/d]"
$ws.Range("D55").Value = "decimal"
$ws.Range("B56").Value = "ZB"
$ws.Range("C56").Value = "fiber intake at baseline [g/d]"
$ws.Range("D56").Value = "decimal"
$ws.Range("B57").Value = "FS"
$ws.Range("C57").Value = "saturated fatty acid intake at baseline [g/d]"
$ws.Range("D57").Value = "decimal"
$ws.Range("B58").Value = "FU"
$ws.Range("C58").Value = "monounsaturated fatty acid intake at baseline [g/d]"
$ws.Range("D58").Value = "decimal"
$ws.Range("B59").Value = "FP"
$ws.Range("C59").Value = "polyunsaturated fatty acid intake at baseline [g/d]"
$ws.Range("D59").Value = "decimal"
$ws.Range("B60").Value = "KD"
$ws.Range("C60").Value = "disaccharide intake at baseline [g/d]"
$ws.Range("D60").Value = "decimal"
$ws.Range("B61").Value = "KM"
$ws.Range("C61").Value = "monosaccharide intake at baseline [g/d]"
$ws.Range("D61").Value = "decimal"
$ws.Range("B62").Value = "KMT"
$ws.Range("C62").Value = "glucose intakeat baseline [g/d]"
$ws.Range("D62").Value = "decimal"
$ws.Range("B63").Value = "KMF"
$ws.Range("C63").Value = "fructose intake at baseline [g/d]"
$ws.Range("D63").Value = "decimal"
$ws.Range("B64").Value = "MNA"
$ws.Range("C64").Value = "sodium intake at baseline [g/d]"
$ws.Range("D64").Value = "decimal"
$ws.Range("B65").Value = "MK"
$ws.Range("C65").Value = "potassium intake at baseline [g/d]"
$ws.Range("D65").Value = "decimal"
$ws.Range("B66").Value = "VEGETABLES_02"
$ws.Range("C66").Value = "Vegetable intake [g/d]"
$ws.Range("D66").Value = "decimal"
$ws.Range("B67").Value = "LEGUMES_TOT_03"
$ws.Range("C67").Value = "Total legumes intake [g/d]"
$ws.Range("D67").Value = "decimal"
$ws.Range("B68").Value = "FRUITS_TOT_04"
$ws.Range("C68").Value = "Total fruit intake [g/d]"
$ws.Range("D68").Value = "decimal"
$ws.Range("B69").Value = "RED_MEAT_0701"
$ws.Range("C69").Value = "Intake of red meat (mammals meat) [g/d]"
$ws.Range("D69").Value = "decimal"
$ws.Range("B70").Value = "PROCMEAT_0704"
$ws.Range("C70").Value = "Intake of processed or preserved meat [g/d]"
$ws.Range("D70").Value = "decimal"
$ws.Range("B71").Value = "SUGAR_CONFECT_11"
$ws.Range("C71").Value = "Intake of sugar and similar, confectionery and water-based sweet desserts [g/d]"
$ws.Range("D71").Value = "decimal"
$ws.Range("B72").Value = "CAKES_12"
$ws.Range("C72").Value = "Intake of cakes and fine bakery products [g/d]"
$ws.Range("D72").Value = "decimal"
$ws.Range("B73").Value = "FRUITVEG_JUICE_1301"
$ws.Range("C73").Value = "Intake of fruit and vegetable juices [g/d]"
$ws.Range("D73").Value = "decimal"
$ws.Range("B74").Value = "SOFTDRINKS_1302"
$ws.Range("C74").Value = "Intake of soft drinks [g/d]"
$ws.Range("D74").Value = "decimal"
$ws.Range("B75").Value = "COFFEE_130301"
$ws.Range("C75").Value = "Coffee intake [g/d]"
$ws.Range("D75").Value = "decimal"
$ws.Range("B76").Value = "TEA_130302"
$ws.Range("C76").Value = "Tea intake [g/d]"
$ws.Range("D76").Value = "decimal"
$ws.Range("B77").Value = "ART_SWEETENER_170201"
$ws.Range("C77").Value = "Intake of artificial sweeteners (e.g., aspartam, saccharine) [g/d]"
$ws.Range("D77").Value = "decimal"

# Remove the trailing row that is now a duplicate of row 77.
$ws.Rows(78).Delete()
